# Auto-generated: apply crypto price/volume updates per diff (Thu Apr 25 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.678.73'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '3.156.04'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = "'612.18"
$ws.Range("E5").Value = '  +1.29%  '
$ws.Range("D6").Value = "'147.98"
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = '3.152.83'
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("D11").Value = "'5.49"
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("D12").Value = "'0.473"
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("E13").Value = '  -2.17%  '
$ws.Range("E14").Value = '  -4.80%  '
$ws.Range("D15").Value = '3.670.32'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("D17").Value = '64.533.77'
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Value = '3.154.15'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = "'6.92"
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = "'8.00"
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").Value = "'13.74"
$ws.Range("E24").Value = '  -3.04%  '
$ws.Range("D25").Value = "'84.09"
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = "'2.83"
$ws.Range("E27").Value = '  -4.39%  '
$ws.Range("D28").Value = "'8.53"
$ws.Range("E28").Value = '  -3.71%  '
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").Value = "'7.03"
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").Value = "'2.09"
$ws.Range("E31").Value = '  -8.48%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").Value = "'26.46"
$ws.Range("E34").Value = '  -2.46%  '
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("D36").Value = '0.0₃0787'
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = "'53.19"
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = "'3.19"
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("D40").Value = "'458.74"
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").Value = "'0.0401"
$ws.Range("E41").Value = '  -1.92%  '
$ws.Range("D42").Value = "'0.120"
$ws.Range("E42").Value = '  -6.86%  '
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").Value = '2.848.72'
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = '  -6.52%  '
$ws.Range("E46").Value = '  -4.05%  '
$ws.Range("E47").Value = '  +3.53%  '
$ws.Range("D48").Value = "'26.52"
$ws.Range("E48").Value = '  -3.98%  '
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").Value = '  -2.20%  '
$ws.Range("D51").Value = "'120.48"
$ws.Range("E51").Value = '  -0.39%  '
